$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "file" primary key / batch date column (column B) values.
# Rows 1-7 previously held "31/05/2023" -> now "19/07/2023"
# Rows 8-11 previously held "01/06/2023" -> now "17/06/2023"
# A leading apostrophe forces the cell to keep its quote-prefixed text
# format (style index 1) instead of reverting to the plain text style.
# (Write the row 8-11 value first so the new shared-string entries are
# appended in the same order as the canonical workbook.)
for ($r = 8; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = "'17/06/2023"
}
for ($r = 1; $r -le 7; $r++) {
    $ws.Cells.Item($r, 2).Value = "'19/07/2023"
}

# Update active selection to match the edited workbook state
$ws.Range("C7").Select() | Out-Null
